# Edit script: append 4 new rows (140-143) of data to worksheet "A"
# matching the commit that extends the data table through row 143.

$data = @(
    @(44579,526,2906,788349,3891821,35368.47,-2380,-215.6919113831161,-96.83737171238673,-118.8545396707294,1101.305078258136,-334.5464510538455,2042.544342360743,-3103472,-43249.70675272169,-19527.63849619965,-23722.06825652203,18231.22934669518,-66971.77500924372,407469.5901211969,36038.72987775276,35895.78566037029,142.9442173824718,36181.67409513523),
    @(44580,1063,2316,1375035,3069343,35028.65,-1253,-319.4227202448045,-154.6455031267674,-164.7772171180371,936.5278611400992,-484.1999373628416,2811.344404997901,-1694308,-208355.5360774495,-103266.6565713897,-105088.8795060599,-86857.65015936468,-313444.4155835094,1788333.174537688,35937.72188997749,35852.42887735178,85.29301262571244,36023.0149026032),
    @(44581,853,2503,1019310,3482133,34715.39,-1650,-452.4804482203241,-229.413227970429,-223.0672202498951,713.4606408902041,-675.5476684702192,3785.796736527683,-2462823,-433802.2824697046,-221244.4737428202,-212557.8087268844,-299415.4588862491,-646360.0911965891,3604796.083341099,35815.48870097974,35795.57693348419,19.9117674955487,35835.40046847529),
    @(44582,680,2757,614754,4860462,34265.37,-2077,-614.9324033982917,-321.7925665719076,-293.1398368263842,420.3208040638199,-908.0722402246759,4954.724496303008,-4245708,-814992.8542227342,-422467.6500556792,-392525.204167055,-691940.6630533042,-1207518.058389789,6642986.024951312,35660.47683088177,35719.06658680997,-58.58975592820934,35601.88707495356),
)


$startRow = 140
$nRows = $data.Count
$nCols = $data[0].Count

# Build a 2D object array for a single bulk write (faster & matches Excel paste semantics)
$arr = New-Object 'object[,]' $nRows, $nCols
for ($i = 0; $i -lt $nRows; $i++) {
    for ($j = 0; $j -lt $nCols; $j++) {
        $arr[$i, $j] = $data[$i][$j]
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A")

$endRow = $startRow + $nRows - 1
$endColLetter = [char](65 + $nCols - 1)
$rangeAddr = "A" + $startRow + ":" + $endColLetter + $endRow

$ws.Range($rangeAddr).Value2 = $arr

# Column A carries the yyyy-mm-dd date number format throughout the table (style index 1);
# apply the same format to the newly added column-A cells.
$dateRangeAddr = "A" + $startRow + ":A" + $endRow
$ws.Range($dateRangeAddr).NumberFormat = "yyyy-mm-dd"
